$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row below row 20 (last data row) for the new period (2508) ---
$ws.Rows("21:21").Insert()

# Copy formatting of (old) row 20 into new row 21 (keeps the "last row" thicker-border style)
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-style row 20 as a "middle" data row (same formatting as rows 16-19) since it is no longer last
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill new row 21 with the same worker data as the rest of the table, for period 2508
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "3885447"
$ws.Range("D21").Value = "PABLO SALGADO ESTRADA"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- 2. Update period labels so they read in ascending order 2503..2507 ---
$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
$ws.Range("E20").Value = "2507"

# --- 3. Swap the "Novedad de Ingreso" / "Novedad de Retiro" column headers ---
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# --- 4. Update the mora value and the period count ---
$ws.Range("E11").Value = 341640
$ws.Range("F13").Value = 6
